# Update the registrant's email address on the "RegistrationData" sheet.
# (cell G1 holds the shared string "rtss59@gmail.com" -> "ghhhh159@gmail.com")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationData")
$ws.Range("G1").Value = "ghhhh159@gmail.com"
